# Fruta / hortaliza, semanal
# Insert a new weekly record as row 2 (shifting the previous rows 2-5 down to 3-6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common values shared by every data row in this subconjunto.
$mercadoId = 11
$mercado   = "Vega Monumental Concepción"
$region    = "Bíobío"
$codreg    = 8
$tipo      = "Fruta"
$productoId = 100107
$producto  = "Otros"
$categoriaId = 100107011
$categoria = "Tuna"
$variedad  = "Sin especificar"
$unidad    = "`$/caja 18 kilos"
$origen    = "Provincia de Melipilla"
$kgUnidad  = 18

# Final data rows 2..6, in order (date, calidad, volumen, precio min, precio max, precio prom, precio/kg)
$rows = @(
    @{ D = 44687; L = "Primera"; M = 100; N = 18000; O = 19000; P = 18500; S = 1028 },
    @{ D = 44280; L = "Primera"; M = 100; N = 14000; O = 15000; P = 14500; S = 806  },
    @{ D = 44280; L = "Segunda"; M = 50;  N = 12000; O = 12000; P = 12000; S = 667  },
    @{ D = 44316; L = "Primera"; M = 50;  N = 20000; O = 20000; P = 20000; S = 1111 },
    @{ D = 44516; L = "Primera"; M = 100; N = 33000; O = 34000; P = 33500; S = 1861 }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $row.D
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $productoId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $kgUnidad
    $r++
}
